$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-11 from 45204 (2023-10-05)
# to 45207 (2023-10-08), keeping existing cell formatting intact.
$ws.Range("C2:C11").Value = 45207
